$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Check In / Check Out Date values (row 2) to new dates, stored as
# text (numFmtId 49 == "@" text format), matching the target workbook.
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "28/05/2016"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "29/05/2016"

# Columns K and L (arrival/depature header + the old duplicated dates) are
# removed entirely from the sheet.
$ws.Range("K1:L2").ClearContents()

# Update the selection to match the saved file (selection was left on the
# now-empty K1:L2 block after the columns were deleted).
$ws.Range("K1:L2").Select()
